$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: lyaktivitet vår, Lars data -- update comment text (style/status unchanged: "påbörjat") ---
$ws.Range("C8").Value = "Saknar vinter/våraktivitet mellan 2011 och 2014. Har bara vinteraktivitet mellan 2000 och 2011 (vet ej vilka exakta datum).  Lars har mailat en del data tidigare år (c:a 2000 - 2005)  till Tomas Meijer och kanske till Anders. Alva kunde inte få ut något vettigt ur databasen. Det fattas årtal och datum på majoriteten av lybesöken. Det jag har nu är vinteraktivitet mellan 2000 och 2010 (BEBODDA_LYOR_HEF 00_10)."

# --- Row 14: GIS-data lyornas avstånd till rödrävslyor -- update comment text (status unchanged: "ej påbörjat") ---
$ws.Range("C14").Value = "det finns rödrävsreproduktion i Peters fil fram till 2008 (röd text) dock är det bara två totalt. Använd koordinaterna för skjutna rävar istället. Finns fram till 2012 i Peters fil.  (Helags_Red_Fox_Feeding).  Maila Lars om resten. Använd Rasmus skript. "

# --- Row 16: Helinventeringar gnagare -- status moves from "ej påbörjat" (red) to "påbörjat" (blue) ---
$ws.Range("B8").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = "påbörjat"
$ws.Range("C16").Value = "Fångstdatan ingår ju i lämmelmodellen så kanske är onödigt? har än så länge bara fångstdata med lyprecision mellan 01 och 04 (Gnagfånst 2001-04 2004-09-16) och 08 och 14 (Sammanfattning08-14.xlsx). Mailat Malin Larm om gnagardata 2015 -2017. Saknas alltså 2005 - 2007. Rasmus tror att datan för de åren kan ha dålig kvalitet."

# --- Row 17: Helinventeringar fågel -- status moves from "ej påbörjat" (red) to "påbörjat" (blue) ---
$ws.Range("B8").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B17").Value = "påbörjat"
$ws.Range("C17").Value = "har mellan 2005 och 2008"

# --- Row 19: Sannolikhet för lämmel inom 1,5 km runt lyan -- status moves from "ej påbörjat" (red) to "klar" (green Brödtext) ---
$ws.Range("B10").Copy()
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("B19").Value = "klar"
$ws.Range("C19").Value = "Klar men pixlarna är endast 49 * 43 m. Inte 500 * 500 meter. 1,5 km eftersom rävar jagar närmare lyan när de har valpar (Frafjord 1993) och måste bära tillbaka mat till lyan (Zapata et al. 1998. Det är vettigare att ta en cirkelradie än riptrianglarna. Dessutom måste jag hålla observationerna oberoende av varandra. Gallant et al (2014) valde max radius på 1,5 km. "

# --- Row 20: Andel bra lämmelhabitat inom området -- status moves from "ej påbörjat" (red) to "påbörjat" (blue); comment unchanged ---
$ws.Range("B8").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value = "påbörjat"

# --- Update the active selection to C16 ---
$ws.Range("C16").Select()
